# The deck has a small "date stamp" textbox on slide 1 (shape "TextBox 4")
# whose second paragraph reads "03-December-24". The edit drops the
# leading "03-" day so it simply reads "December-24".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(2)
$run = $para.Runs(1)

if ($run.Text -eq "03-December-24") {
    $run.Text = "December-24"
} else {
    # Fall back to a targeted find/replace in case the run layout differs.
    $found = $tr.Find("03-December-24")
    if ($found -ne $null) {
        $found.Text = "December-24"
    }
}
